# Criação de todas sheets que precisam de testes
# Renomeia a planilha original e cria as demais planilhas (sheets) do fluxo de UAT.

$wb = $excel.ActiveWorkbook

# 1) Renomeia a primeira planilha (antiga "Planilha1") para "Cadastro"
$cadastro = $wb.Worksheets.Item(1)
$cadastro.Name = "Cadastro"

# 2) Cria as novas planilhas, na ordem em que foram criadas originalmente,
#    sempre inserindo após a última planilha existente.
$novosNomes = @(
    "Recuperação de Senha",
    "Login",
    "Home",
    "Configurações",
    "Perfil do Jogador",
    "Perfil da Equipe",
    "Continuação do Cadastro",
    "Buscar",
    "Agendamento de Partida"
)

$ultima = $wb.Worksheets.Item($wb.Worksheets.Count)
foreach ($nome in $novosNomes) {
    $nova = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ultima)
    $nova.Name = $nome
    $ultima = $nova
}

# 3) Move "Continuação do Cadastro" para logo após "Cadastro"
$continuacao = $wb.Worksheets.Item("Continuação do Cadastro")
$continuacao.Move([System.Reflection.Missing]::Value, $cadastro)

# 4) Ajusta a seleção da planilha "Login" (equivalente ao bloco copiado de Cadastro)
$login = $wb.Worksheets.Item("Login")
$login.Range("A1:C8").Select() | Out-Null

# 5) Reativa "Cadastro" como planilha ativa e atualiza a célula selecionada
$cadastro.Activate() | Out-Null
$cadastro.Range("C15").Select() | Out-Null
